$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1627906976744186
$ws.Range("C2").Value = 0.6085271317829457
$ws.Range("J2").Value = 0.02713178294573643
$ws.Range("P2").Value = 0.1085271317829457
$ws.Range("S2").Value = 0.09302325581395349

$ws.Range("B3").Value = 0.01818181818181818
$ws.Range("C3").Value = 0.0303030303030303
$ws.Range("J3").Value = 0.02424242424242424
$ws.Range("P3").Value = 0.7939393939393939
$ws.Range("S3").Value = 0.1333333333333333

$ws.Range("B6").Value = 0.04587155963302753
$ws.Range("D6").Value = 0.009174311926605505
$ws.Range("F6").Value = 0.05045871559633028
$ws.Range("J6").Value = 0.2706422018348624
$ws.Range("O6").Value = 0.01834862385321101
$ws.Range("Q6").Value = 0.1422018348623853
$ws.Range("R6").Value = 0.05504587155963303
$ws.Range("S6").Value = 0.4082568807339449

$ws.Range("B7").Value = 0.09950248756218906
$ws.Range("D7").Value = 0.01990049751243781
$ws.Range("E7").Value = 0.004975124378109453
$ws.Range("F7").Value = 0.05970149253731343
$ws.Range("J7").Value = 0.1691542288557214
$ws.Range("O7").Value = 0.02985074626865672
$ws.Range("Q7").Value = 0.1293532338308458
$ws.Range("R7").Value = 0.09950248756218906
$ws.Range("S7").Value = 0.3880597014925373

$ws.Range("B8").Value = 0.08078602620087336
$ws.Range("D8").Value = 0.01965065502183406
$ws.Range("F8").Value = 0.07423580786026202
$ws.Range("J8").Value = 0.1179039301310044
$ws.Range("O8").Value = 0.01746724890829694
$ws.Range("Q8").Value = 0.1877729257641921
$ws.Range("R8").Value = 0.08078602620087336
$ws.Range("S8").Value = 0.4213973799126637

$ws.Range("B9").Value = 0.08187134502923976
$ws.Range("D9").Value = 0.01754385964912281
$ws.Range("F9").Value = 0.05847953216374269
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.03508771929824561
$ws.Range("Q9").Value = 0.1988304093567251
$ws.Range("R9").Value = 0.0935672514619883
$ws.Range("S9").Value = 0.4035087719298245

$ws.Range("B10").Value = 0.1111111111111111
$ws.Range("D10").Value = 0.02136752136752137
$ws.Range("F10").Value = 0.08632478632478632
$ws.Range("J10").Value = 0.1222222222222222
$ws.Range("O10").Value = 0.0188034188034188
$ws.Range("Q10").Value = 0.1974358974358974
$ws.Range("R10").Value = 0.07777777777777778
$ws.Range("S10").Value = 0.3649572649572649

$ws.Range("G11").Value = 0.1423487544483986
$ws.Range("J11").Value = 0.099644128113879
$ws.Range("K11").Value = 0.2135231316725979
$ws.Range("L11").Value = 0.5338078291814946
$ws.Range("S11").Value = 0.01067615658362989

$ws.Range("G12").Value = 0.7870967741935484
$ws.Range("J12").Value = 0.1806451612903226
$ws.Range("L12").Value = 0.01935483870967742
$ws.Range("S12").Value = 0.01290322580645161

$ws.Range("G13").Value = 0.7586206896551724
$ws.Range("J13").Value = 0.2241379310344828
$ws.Range("S13").Value = 0.01724137931034483

$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.3333333333333333

$ws.Range("F15").Value = 0.004524886877828055
$ws.Range("H15").Value = 0.16289592760181
$ws.Range("I15").Value = 0.04977375565610859
$ws.Range("J15").Value = 0.3981900452488688
$ws.Range("K15").Value = 0.04524886877828054
$ws.Range("M15").Value = 0.004524886877828055
$ws.Range("N15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.04072398190045249
$ws.Range("S15").Value = 0.2895927601809955

$ws.Range("F16").Value = 0.03296703296703297
$ws.Range("H16").Value = 0.2032967032967033
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.3681318681318682
$ws.Range("K16").Value = 0.0989010989010989
$ws.Range("M16").Value = 0.04945054945054945
$ws.Range("O16").Value = 0.05494505494505494
$ws.Range("S16").Value = 0.1208791208791209

$ws.Range("F17").Value = 0.01234567901234568
$ws.Range("H17").Value = 0.1703703703703704
$ws.Range("I17").Value = 0.1012345679012346
$ws.Range("J17").Value = 0.4271604938271605
$ws.Range("K17").Value = 0.08888888888888889
$ws.Range("M17").Value = 0.01728395061728395
$ws.Range("O17").Value = 0.05925925925925926
$ws.Range("S17").Value = 0.1234567901234568

$ws.Range("F18").Value = 0.01123595505617977
$ws.Range("H18").Value = 0.2191011235955056
$ws.Range("I18").Value = 0.07865168539325842
$ws.Range("J18").Value = 0.398876404494382
$ws.Range("K18").Value = 0.1067415730337079
$ws.Range("M18").Value = 0.005617977528089887
$ws.Range("N18").Value = 0.005617977528089887
$ws.Range("O18").Value = 0.07865168539325842
$ws.Range("S18").Value = 0.09550561797752809

$ws.Range("F19").Value = 0.01483924154987634
$ws.Range("H19").Value = 0.2308326463314097
$ws.Range("I19").Value = 0.07749381698268755
$ws.Range("J19").Value = 0.3322341302555647
$ws.Range("K19").Value = 0.1129431162407255
$ws.Range("M19").Value = 0.0313272877164056
$ws.Range("N19").Value = 0.0008244023083264633
$ws.Range("O19").Value = 0.07007419620774938
$ws.Range("S19").Value = 0.1294311624072547
